# Update cryptocurrency price/volume data as reflected in the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.016.69'
$ws.Range('E2').Value = '  -2.12%  '
$ws.Range('D3').Value = '1.666.68'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '216.52'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.63%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5100'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.59%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.004'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2639'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.66%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06393'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.71%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.91'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.40%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07397'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.57%  '
$ws.Range('D12').Value = '1.679.44'
$ws.Range('E12').Value = '  -1.15%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.499'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E14').Value = '  -0.44%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000008483'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.35%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.16'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.34%  '
$ws.Range('D17').Value = '26.055.06'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '4.927'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.07%  '
$ws.Range('E19').Value = '  -0.07%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.73'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.34%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '189.58'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.21%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.206'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.25%  '
$ws.Range('E23').Value = '  -0.15%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '145.00'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.11%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '7.603'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.47%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1193'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.36%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.63'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.52%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06645'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +16.28%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.313'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.19%  '
$ws.Range('E30').Value = '  -2.03%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.519'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.36%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.507'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.06%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.634'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.00%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.016'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.40%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6079'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.72%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.367'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.25%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.709'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.97%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.207'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +6.13%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01606'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.26%  '
$ws.Range('D40').Value = '1.075.61'
$ws.Range('E40').Value = '  -2.49%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8575'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.45%  '
$ws.Range('E42').Value = '  +0.28%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '100.44'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').Value = '1.814.21'
$ws.Range('E44').Value = '  -2.36%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00000000112'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.69%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '56.27'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.85%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.008'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.35%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.022'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.95%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05207'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.78%  '
$ws.Range('E50').Value = '  -0.91%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '5.950'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.52%  '
